$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet2")
$ws.Activate()

# Add new meeting-diary entry as row 11, copying the formatting of the most
# similarly-formatted existing entry (row 9) so styles are reused instead of
# duplicated.
$ws.Range("A9:E9").Copy()
$ws.Range("A11:E11").PasteSpecial(-4122)

$ws.Range("A11").Value = 45225
$ws.Range("B11").Value = 0.39583333333333331
$ws.Range("C11").Value = 0.45833333333333331
$ws.Range("D11").Value = "All"
$ws.Range("E11").Value = "Preparing for the presentation and checking the final report."

# Widen column E to fit the new, longer discussion text.
$ws.Columns("E").ColumnWidth = 50

# Update the saved cursor/selection position.
$ws.Range("E22").Select()

Write-Host "done"
